$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct two daily scores on 2025-02-12
$ws.Range("K35").Value = $false
$ws.Range("K36").Value = $true

# Append new daily scores for 2025-02-13
$ws.Range("A38:A40").NumberFormat = "@"
$ws.Range("A38").Value = "2025-02-13"
$ws.Range("B38").Value = "sleep"
$ws.Range("C38").Value = $false
$ws.Range("D38").Value = $true
$ws.Range("E38").Value = $true
$ws.Range("F38").Value = $false
$ws.Range("G38").Value = $true
$ws.Range("H38").Value = $false
$ws.Range("I38").Value = $true
$ws.Range("J38").Value = $true
$ws.Range("K38").Value = $true
$ws.Range("L38").Value = $false
$ws.Range("M38").Value = $true
$ws.Range("N38").Value = $true
$ws.Range("O38").Value = $true

$ws.Range("A39").Value = "2025-02-13"
$ws.Range("B39").Value = "activity"
$ws.Range("C39").Value = $false
$ws.Range("D39").Value = $false
$ws.Range("E39").Value = $true
$ws.Range("F39").Value = $true
$ws.Range("G39").Value = $true
$ws.Range("H39").Value = $true
$ws.Range("I39").Value = $true
$ws.Range("J39").Value = $false
$ws.Range("K39").Value = $false
$ws.Range("L39").Value = $true
$ws.Range("M39").Value = $false
$ws.Range("N39").Value = $false
$ws.Range("O39").Value = $false

$ws.Range("A40").Value = "2025-02-13"
$ws.Range("B40").Value = "weekly_activity"
$ws.Range("C40").Value = $false
$ws.Range("D40").Value = $false
$ws.Range("E40").Value = $true
$ws.Range("F40").Value = $true
$ws.Range("G40").Value = $true
$ws.Range("H40").Value = $false
$ws.Range("I40").Value = $true
$ws.Range("J40").Value = $true
$ws.Range("K40").Value = $true
$ws.Range("L40").Value = $true
$ws.Range("M40").Value = $true
$ws.Range("N40").Value = $false
$ws.Range("O40").Value = $false

# Restore default styling on the date column (avoid leftover date/text format)
$ws.Range("A38:A40").Style = "Normal"
